$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per the target diff (rows 2-25 of the rate table)
$ws.Range("B2").Value = 16500
$ws.Range("S2").Value = 36728
$ws.Range("T2").Value = 5799.800999999999
$ws.Range("U2").Value = 66500
$ws.Range("B3").Value = 13500
$ws.Range("I3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 35269
$ws.Range("T3").Value = 5649.136500000001
$ws.Range("U3").Value = 63500
$ws.Range("X3").Value = 10000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 33348
$ws.Range("T4").Value = 5554.8325
$ws.Range("U4").Value = 42500
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 31700
$ws.Range("T5").Value = 5481.98
$ws.Range("U5").Value = 42500
$ws.Range("S6").Value = 30762
$ws.Range("T6").Value = 5456.8815
$ws.Range("P7").Value = 0
$ws.Range("S7").Value = 30953
$ws.Range("T7").Value = 5512.395
$ws.Range("U7").Value = 22500
$ws.Range("P8").Value = 0
$ws.Range("S8").Value = 31330
$ws.Range("T8").Value = 5765.477227722772
$ws.Range("U8").Value = 22500
$ws.Range("P9").Value = 0
$ws.Range("S9").Value = 32595
$ws.Range("T9").Value = 6519.4465
$ws.Range("U9").Value = 22500
$ws.Range("P10").Value = 0
$ws.Range("S10").Value = 36208
$ws.Range("T10").Value = 7901.5895
$ws.Range("U10").Value = 22500
$ws.Range("P11").Value = 10000
$ws.Range("S11").Value = 39069
$ws.Range("T11").Value = 13052.0005
$ws.Range("U11").Value = 32500
$ws.Range("P12").Value = 10000
$ws.Range("S12").Value = 39884
$ws.Range("T12").Value = 14886.7985
$ws.Range("U12").Value = 32500
$ws.Range("P13").Value = 10000
$ws.Range("S13").Value = 39218
$ws.Range("T13").Value = 14810.005
$ws.Range("U13").Value = 32500
$ws.Range("P14").Value = 10000
$ws.Range("S14").Value = 38683
$ws.Range("T14").Value = 14699.3385
$ws.Range("U14").Value = 32500
$ws.Range("B15").Value = 12500
$ws.Range("I15").Value = 10000
$ws.Range("S15").Value = 40878
$ws.Range("T15").Value = 15357.636
$ws.Range("U15").Value = 52500
$ws.Range("S16").Value = 41508
$ws.Range("T16").Value = 15432.207
$ws.Range("S17").Value = 40766
$ws.Range("T17").Value = 15338.9075
$ws.Range("S18").Value = 38994
$ws.Range("T18").Value = 15792.80217625723
$ws.Range("S19").Value = 36246
$ws.Range("T19").Value = 15557.50171551809
$ws.Range("V19").Value = 8357.577825000002
$ws.Range("W19").Value = -2.485794135445751
$ws.Range("T20").Value = 15152.49877462994
$ws.Range("T21").Value = 13812.90902852661
$ws.Range("T22").Value = 12007.84432898735
$ws.Range("T23").Value = 9523.143
$ws.Range("B24").Value = 22000
$ws.Range("T24").Value = 6375.7855
$ws.Range("U24").Value = 72000
$ws.Range("B25").Value = 12500
$ws.Range("T25").Value = 5494.996500000001
$ws.Range("U25").Value = 52500

# Row 15: remove V15 and W15 (WESM_RATE / CURRENT_RATE no longer computed for this hour)
$ws.Range("V15").ClearContents()
$ws.Range("W15").ClearContents()

